$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.105.19"
$ws.Range("E2").Value = "  +3.33%  "
$ws.Range("D3").Value = "2.457.18"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.74"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.24"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.10%  "
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +5.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.16"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.43"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.08"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "2.840.75"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "2.463.04"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "46.045.80"
$ws.Range("E18").Value = "  +3.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.60"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "0.0₃0933"
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.67"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.93%  "
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "248.15"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.63%  "
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.02"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.87%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  -4.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.70"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.92"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "49.27"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("E32").Value = "  +5.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.41"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.66%  "
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.91"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "127.36"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.28"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.04"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0292"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("D45").Value = "1.962.66"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("E47").Value = "  -2.73%  "
$ws.Range("E48").Value = "  +10.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.23"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "77.75"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.05%  "
$ws.Range("E51").Value = "  +5.82%  "
